$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: current_phase (D33) bumped from 1 to 2
$ws.Range("D33").Value = 2

# Row 34: bot advanced this group's phase, recorded a new last_action_date,
# and counted one reply (message id 103)
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = "2026-02-12T22:38:11.289567+00:00"
$ws.Range("I34").Value = 1
$ws.Range("M34").Value = "[103]"
